$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grades")

# Quizzes block - HW-M8 (row 13) got graded: 20/20
$ws.Range("G13").Value = 20

# Quiz Chapter 7 (row 13) and Quiz Chapter 8 (row 14) got graded: 13/14 and 13/13
$ws.Range("L13").Value = 13
$ws.Range("L14").Value = 13

# Assignments block - remaining discussion rows got graded: 10/10 each
$ws.Range("B15").Value = 10
$ws.Range("B16").Value = 10
$ws.Range("B24").Value = 10
$ws.Range("B25").Value = 10
$ws.Range("B27").Value = 10

# Quizzes block - Quiz Chapter 5 (row 15) got graded: 20/20
$ws.Range("G15").Value = 20

# Update the active selection to V7 as reflected in the saved view
$ws.Range("V7").Select()
